$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (W1) onto the
# new header cell (X1) so it keeps the bold/border/centered style without
# allocating a brand-new style entry.
$ws.Range("W1").Copy()
$ws.Range("X1").PasteSpecial(-4122)
$ws.Range("X1").Value = "l1"

# Fill the new column with 0 for every data row (rows 2-77), matching the
# existing V/W "numeric flag" columns.
for ($r = 2; $r -le 77; $r++) {
    $ws.Cells.Item($r, 24).Value = 0
}
